$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Simple value corrections (rows before the inserted row, so no shift needed)
$ws1.Range("M218").Value = 829.4400000000001
$ws1.Range("M231").Value = 6635.52
$ws1.Range("L302").Value = 179.12

# A new salesperson/client combo row was inserted right before "VACA PANCHI
# DORYS CAROLINA" (old row 309), pushing everything from row 309 down by one.
$ws1.Rows.Item(309).Insert()

$ws1.Range("A309").Value = "OFICINA-CATAECSA"
$ws1.Range("B309").Value = "TORRES CADENA JAVIER JOSUE"
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(309, $c).Value = 0
}

# ---------------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Simple value corrections (rows before the inserted row, so no shift needed)
$ws2.Range("F222").Value = 829.4400000000001
$ws2.Range("F235").Value = 6635.52
$ws2.Range("F306").Value = 179.12

# Same new row, inserted right before "VACA PANCHI DORYS CAROLINA" (old row 313)
$ws2.Rows.Item(313).Insert()

$ws2.Range("A313").Value = "OFICINA-CATAECSA"
$ws2.Range("B313").Value = "TORRES CADENA JAVIER JOSUE"
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(313, $c).Value = 0
}

# ---------------------------------------------------------------------------
# Sheet 3: "CUMPLIMIENTO MENSUAL"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D60").Value = 7464.96
$ws3.Range("E60").Value = 31077.29
$ws3.Range("F60").Value = 0.193682517237577

$ws3.Range("D62").Value = 2818.29
$ws3.Range("E62").Value = -2818.29

$ws3.Range("D76").Value = 18778.57
$ws3.Range("E76").Value = 388833.3870193433
$ws3.Range("F76").Value = 0.04606972311930697
